$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "29.197.81"
$ws.Range("E2").Value = "  +0.11%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "1.852.30"
$ws.Range("E3").Value = "  -0.13%  "

$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "1.000"
$ws.Range("E4").Value = "  +0.01%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "0.6992"
$ws.Range("E5").Value = "  +1.71%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "237.43"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("E7").Value = "  -0.02%  "

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07877"
$ws.Range("E8").Value = "  +1.25%  "

$ws.Range("E9").Value = "  -0.87%  "

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "23.91"
$ws.Range("E10").Value = "  +3.06%  "

$ws.Range("E11").Value = "  +0.60%  "

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "1.848.82"
$ws.Range("E12").Value = "  +0.78%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "5.186"
$ws.Range("E13").Value = "  -0.25%  "

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "0.7058"
$ws.Range("E14").Value = "  -2.19%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "89.34"
$ws.Range("E15").Value = "  +0.04%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "29.234.30"
$ws.Range("E16").Value = "  +0.23%  "

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "5.803"
$ws.Range("E17").Value = "  +1.15%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "0.000007828"
$ws.Range("E18").Value = "  +0.39%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "13.20"
$ws.Range("E19").Value = "  -0.68%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "235.44"
$ws.Range("E20").Value = "  +0.62%  "

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9998"
$ws.Range("E21").Value = "  -0.02%  "

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "2.104.39"
$ws.Range("E22").Value = "  -0.20%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "1.000"
$ws.Range("E23").Value = "  -0.02%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "7.487"
$ws.Range("E24").Value = "  +0.16%  "

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "162.51"
$ws.Range("E25").Value = "  +0.42%  "

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "8.886"
$ws.Range("E26").Value = "  -0.99%  "

$ws.Range("E27").Value = "  -0.77%  "

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "18.02"
$ws.Range("E28").Value = "  -0.05%  "

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "1.914"
$ws.Range("E29").Value = "  -2.23%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "1.400"
$ws.Range("E30").Value = "  -0.52%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "1.474"
$ws.Range("E31").Value = "  -0.43%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "4.297"
$ws.Range("E32").Value = "  -4.62%  "

$ws.Range("E33").Value = "  +0.01%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "0.05144"
$ws.Range("E34").Value = "  -1.33%  "

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "1.165"
$ws.Range("E35").Value = "  -1.11%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "0.7065"
$ws.Range("E36").Value = "  +0.33%  "

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9966"
$ws.Range("E37").Value = "  -0.40%  "

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "2.676"
$ws.Range("E38").Value = "  +0.23%  "

$ws.Range("E39").Value = "  -0.17%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "2.703"
$ws.Range("E40").Value = "  +0.38%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "1.150.34"
$ws.Range("E41").Value = "  +4.45%  "

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9202"
$ws.Range("E42").Value = "  -1.73%  "

$ws.Range("E43").Value = "  +0.74%  "

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "0.4235"
$ws.Range("E44").Value = "  -1.12%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "69.95"
$ws.Range("E45").Value = "  -0.52%  "

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9997"
$ws.Range("E46").Value = "  -0.05%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "102.88"
$ws.Range("E47").Value = "  +0.45%  "

$ws.Range("E48").Value = "  -2.87%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "1.735"
$ws.Range("E49").Value = "  -3.35%  "

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "9.158"
$ws.Range("E50").Value = "  +0.00%  "

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "6.954"
$ws.Range("E51").Value = "  -0.58%  "
